# Applies the diff:
#  - The paragraph "Statusbericht Changes aufführen -> Wir vermuten das..."
#    is split into three runs, wrapping "Changes" in spell-check proofErr
#    tags (the visible text content itself is unchanged).
#  - The last paragraph (previously empty, holding only the _GoBack
#    bookmark) gets new text:
#    "Statusbericht fordert im Normalfall genehmigung für weiteres
#    vorgehen..." - split into runs around the existing _GoBack bookmark
#    and around the words "genehmigung"/"vorgehen", which get wrapped in
#    proofErr tags.
#
# Both paragraphs are rewritten wholesale via Range.InsertXML so the
# exact run/proofErr/bookmark structure from the target revision can be
# produced.

$d = $word.ActiveDocument

$wNs   = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# Locate the two paragraphs we need to touch by their current content,
# rather than relying on a fixed index.
$changesParaIndex = -1
$bookmarkParaIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Statusbericht Changes*") {
        $changesParaIndex = $i
    }
    if ($t.Trim() -eq "") {
        $bookmarkParaIndex = $i
    }
}

# --- "Statusbericht Changes aufführen -> Wir vermuten das..." paragraph ---
$pChanges = $d.Paragraphs($changesParaIndex)
$rChanges = $pChanges.Range
$attrsChanges = 'w14:paraId="4046A66B" w14:textId="77777777" w:rsidR="00B67CE7" w:rsidRDefault="00B67CE7" w:rsidP="00C10240"'
$bodyChanges = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Statusbericht </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Changes</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> auff&#252;hren -&gt; Wir vermuten das&#8230;</w:t></w:r>'
$xmlChanges = "<w:p $wNs $w14Ns $attrsChanges>$bodyChanges</w:p>"
$rChanges.InsertXML($xmlChanges)

# --- Empty paragraph that only held the _GoBack bookmark ---
$pBookmark = $d.Paragraphs($bookmarkParaIndex)
$rBookmark = $pBookmark.Range
$attrsBookmark = 'w14:paraId="4F94EC54" w14:textId="77777777" w:rsidR="007A253A" w:rsidRDefault="007A253A" w:rsidP="00C10240"'
$bodyBookmark = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>Statusbericht fordert im Normalf</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t xml:space="preserve">all </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>genehmigung</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> f&#252;r weiteres </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>vorgehen</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>&#8230;</w:t></w:r>'
$xmlBookmark = "<w:p $wNs $w14Ns $attrsBookmark>$bodyBookmark</w:p>"
$rBookmark.InsertXML($xmlBookmark)
